# Katalon IWP EmulatorData.xlsx update
# Commit: "Added ABP test cases and modified IWP Bootstrap deferred test cases"
#
# Content change: a new AutoPay-All-Fields test-case row is appended
# (row 34) below the existing data (rows 1-33), cloning the shape of the
# existing "AutoPay All Fields" row (row 4) but with a new ID ("33") and
# a new MV ("2.5"). The view/selection is also nudged to reflect where
# the author was working when the row was added.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 34: AutoPay All Fields test case (ID 33, MV 2.5) ---------
$ws.Range("A34").Value2 = "AutoPay All Fields"
$ws.Range("C34").Value2 = "33"
$ws.Range("D34").Value2 = "2.5"
$ws.Range("E34").Value2 = "13.50"
$ws.Range("G34").Value2 = "AutoPay"
$ws.Range("H34").Value2 = "en_US"
$ws.Range("I34").Value2 = "Jonty"
$ws.Range("J34").Value2 = "Smith"
$ws.Range("K34").Value2 = "15 Elm St"
$ws.Range("L34").Value2 = "Suite 600"
$ws.Range("M34").Value2 = "840"
$ws.Range("N34").Value2 = "Gambrills"
$ws.Range("O34").Value2 = "MD"
$ws.Range("P34").Value2 = "21054"
$ws.Range("S34").Value2 = "iahmed@govolution.com"
$ws.Range("T34").Value2 = "udf data 1"
$ws.Range("U34").Value2 = "udf data 2"
$ws.Range("V34").Value2 = "udf data 3"
$ws.Range("W34").Value2 = "udf data 4"
$ws.Range("X34").Value2 = "udf data 5"
$ws.Range("Y34").Value2 = "udf data 6"
$ws.Range("Z34").Value2 = "Orange"
$ws.Range("AA34").Value2 = "Soccer"
$ws.Range("AB34").Value2 = "udf data 9"
$ws.Range("AC34").Value2 = "udf data 10"

# S34 keeps the same bordered/wrap-text style used by every other cell in
# the "Email" (S) column (the plain Value2 write above only inherits the
# column's default style, so the formatting is copied across explicitly).
$ws.Range("S4").Copy() | Out-Null
$ws.Range("S34").PasteSpecial(-4122) | Out-Null

# --- View / selection: scrolled down to the new row, cell J39 active --
$ws.Range("A27").Select() | Out-Null
$ws.Range("J39").Select() | Out-Null
